$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows there currently are (column A has the file names).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add a new "Electrode Locations" header in column C, matching the bold/centered/
# bordered formatting already used for the "File Name" / "Unnormalized P_max" headers.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Electrode Locations"

# Populate the new Electrode Locations column from each row's file name
# (the electrode location is the text before the first underscore, e.g.
# "A11_monopolar_10V_5Hz.txt" -> "A11").
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = $ws.Cells.Item($r, 1).Value2
    $location = $fileName.Split("_")[0]
    $ws.Cells.Item($r, 3).Value = $location
}

# Build a hidden helper column (F) with a single numeric sort key that orders the
# electrode locations the natural way (A1, A5, A7, ..., A13, B3, B15, C1, ... O13)
# instead of a plain alphabetical text sort (which would incorrectly put "A11"
# before "A5"). The key combines the letter (its character code) and the numeric
# part of the location, e.g. "A11" -> 65*100 + 11 = 6511.
for ($r = 2; $r -le $lastRow; $r++) {
    $location = $ws.Cells.Item($r, 3).Value2
    $letter = $location.Substring(0, 1)
    $number = $location.Substring(1)
    $letterCode = [int][char]$letter
    $sortKey = ($letterCode * 100) + [int]$number
    $ws.Cells.Item($r, 6).Value = $sortKey
}

# Sort all of the data rows (file name, value, electrode location, helper key)
# in ascending order of the helper key, i.e. by electrode location from A1 to O15.
$sortRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 6))
$sortKeyRange = $ws.Range($ws.Cells.Item(2, 6), $ws.Cells.Item($lastRow, 6))
$sortRange.Sort($sortKeyRange, 1)

# Remove the temporary helper column now that the rows are in the right order.
$ws.Range($ws.Cells.Item(1, 6), $ws.Cells.Item($lastRow, 6)).Clear()
